$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4694920778274536
$ws.Range("B1").Value = 2.663505792617798
$ws.Range("C1").Value = 6.136752128601074
$ws.Range("D1").Value = 1.734527945518494
$ws.Range("E1").Value = 1.027464151382446
